$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "SUBSET"
$ws.Range("F2").Value = 0.8769230769230763
$ws.Range("G2").Value = 0.9548980933596314
$ws.Range("C3").Value = "SUBSET"
$ws.Range("F3").Value = 0.8615384615384606
$ws.Range("G3").Value = 0.9539776462853382
$ws.Range("C4").Value = "SUBSET"
$ws.Range("F4").Value = 0.8743589743589736
$ws.Range("G4").Value = 0.97120315581854
$ws.Range("H4").Value = "ENSG00000203875, ENSG00000269893, ENSG00000204387, ENSG00000233493, ENSG00000215908"
$ws.Range("I4").Value = "SNHG5, SNHG8, C6orf48, TMEM238, CROCCP2"
$ws.Range("C5").Value = "SUBSET"
$ws.Range("F5").Value = 0.8923076923076918
$ws.Range("G5").Value = 0.9697567389875081
$ws.Range("H5").Value = "ENSG00000203875, ENSG00000269893, ENSG00000274012, ENSG00000258920, ENSG00000204387, ENSG00000215908"
$ws.Range("I5").Value = "SNHG5, SNHG8, RN7SL2, FOXN3-AS1, C6orf48, CROCCP2"
$ws.Range("C6").Value = "SUBSET"
$ws.Range("F6").Value = 0.8999999999999997
$ws.Range("G6").Value = 0.9760683760683755
$ws.Range("H6").Value = "ENSG00000203875, ENSG00000269893, ENSG00000274012, ENSG00000258920, ENSG00000204387, ENSG00000233493, ENSG00000215908"
$ws.Range("I6").Value = "SNHG5, SNHG8, RN7SL2, FOXN3-AS1, C6orf48, TMEM238, CROCCP2"
$ws.Range("C7").Value = "SUBSET"
$ws.Range("F7").Value = 0.871794871794871
$ws.Range("G7").Value = 0.9584483892176197
$ws.Range("H7").Value = "ENSG00000203875, ENSG00000269893, ENSG00000236552, ENSG00000258920, ENSG00000204387, ENSG00000215908, ENSG00000234741, ENSG00000179085"
$ws.Range("I7").Value = "SNHG5, SNHG8, RPL13AP5, FOXN3-AS1, C6orf48, CROCCP2, GAS5, DPM3"
$ws.Range("C8").Value = "SUBSET"
$ws.Range("F8").Value = 0.8846153846153839
$ws.Range("G8").Value = 0.9742274819197891
$ws.Range("H8").Value = "ENSG00000203875, ENSG00000269893, ENSG00000236552, ENSG00000274012, ENSG00000255559, ENSG00000204387, ENSG00000233493, ENSG00000215908, ENSG00000234741"
$ws.Range("I8").Value = "SNHG5, SNHG8, RPL13AP5, RN7SL2, ZNF252P-AS1, C6orf48, TMEM238, CROCCP2, GAS5"
$ws.Range("C9").Value = "SUBSET"
$ws.Range("F9").Value = 0.8871794871794867
$ws.Range("G9").Value = 0.9733070348454961
$ws.Range("H9").Value = "ENSG00000203875, ENSG00000269893, ENSG00000274012, ENSG00000258920, ENSG00000225864, ENSG00000204387, ENSG00000233493, ENSG00000215908, ENSG00000234741, ENSG00000179085"
$ws.Range("I9").Value = "SNHG5, SNHG8, RN7SL2, FOXN3-AS1, HCG4P11, C6orf48, TMEM238, CROCCP2, GAS5, DPM3"
$ws.Range("C10").Value = "SUBSET"
$ws.Range("F10").Value = 0.8846153846153839
$ws.Range("G10").Value = 0.9642340565417487
$ws.Range("H10").Value = "ENSG00000203875, ENSG00000269893, ENSG00000236552, ENSG00000274012, ENSG00000278771, ENSG00000272906, ENSG00000226287, ENSG00000204387, ENSG00000233493, ENSG00000215908, ENSG00000179085"
$ws.Range("I10").Value = "SNHG5, SNHG8, RPL13AP5, RN7SL2, Metazoa_SRP, RP11-533E19.7, TMEM191A, C6orf48, TMEM238, CROCCP2, DPM3"
$ws.Range("C11").Value = "SUBSET"
$ws.Range("F11").Value = 0.8794871794871788
$ws.Range("G11").Value = 0.9654832347140034
$ws.Range("H11").Value = "ENSG00000203875, ENSG00000269893, ENSG00000236552, ENSG00000274012, ENSG00000272906, ENSG00000258920, ENSG00000226287, ENSG00000204387, ENSG00000233493, ENSG00000215908, ENSG00000215414, ENSG00000230979"
$ws.Range("I11").Value = "SNHG5, SNHG8, RPL13AP5, RN7SL2, RP11-533E19.7, FOXN3-AS1, TMEM191A, C6orf48, TMEM238, CROCCP2, PSMA6P1, AC079250.1"
$ws.Range("C12").Value = "SUBSET"
$ws.Range("F12").Value = 0.8846153846153839
$ws.Range("G12").Value = 0.9679158448389215
$ws.Range("H12").Value = "ENSG00000203875, ENSG00000269893, ENSG00000236552, ENSG00000274012, ENSG00000278771, ENSG00000272906, ENSG00000255559, ENSG00000226287, ENSG00000204387, ENSG00000233493, ENSG00000215908, ENSG00000234741, ENSG00000179085"
$ws.Range("I12").Value = "SNHG5, SNHG8, RPL13AP5, RN7SL2, Metazoa_SRP, RP11-533E19.7, ZNF252P-AS1, TMEM191A, C6orf48, TMEM238, CROCCP2, GAS5, DPM3"
$ws.Range("C13").Value = "SUBSET"
$ws.Range("F13").Value = 0.8871794871794865
$ws.Range("G13").Value = 0.9671268902038129
$ws.Range("H13").Value = "ENSG00000203875, ENSG00000269893, ENSG00000236552, ENSG00000274012, ENSG00000278771, ENSG00000255559, ENSG00000258920, ENSG00000226287, ENSG00000204387, ENSG00000233493, ENSG00000215908, ENSG00000234741, ENSG00000230979, ENSG00000179085"
$ws.Range("I13").Value = "SNHG5, SNHG8, RPL13AP5, RN7SL2, Metazoa_SRP, ZNF252P-AS1, FOXN3-AS1, TMEM191A, C6orf48, TMEM238, CROCCP2, GAS5, AC079250.1, DPM3"
$ws.Range("C14").Value = "SUBSET"
$ws.Range("F14").Value = 0.8871794871794865
$ws.Range("G14").Value = 0.9708086785009856
$ws.Range("H14").Value = "ENSG00000203875, ENSG00000269893, ENSG00000236552, ENSG00000274012, ENSG00000278771, ENSG00000272906, ENSG00000255559, ENSG00000258920, ENSG00000225864, ENSG00000226287, ENSG00000204387, ENSG00000233493, ENSG00000215908, ENSG00000234741, ENSG00000179085"
$ws.Range("I14").Value = "SNHG5, SNHG8, RPL13AP5, RN7SL2, Metazoa_SRP, RP11-533E19.7, ZNF252P-AS1, FOXN3-AS1, HCG4P11, TMEM191A, C6orf48, TMEM238, CROCCP2, GAS5, DPM3"
$ws.Range("C15").Value = "SUBSET"
$ws.Range("G15").Value = 0.9693622616699532
$ws.Range("H15").Value = "ENSG00000203875, ENSG00000269893, ENSG00000236552, ENSG00000274012, ENSG00000272906, ENSG00000255559, ENSG00000258920, ENSG00000226287, ENSG00000204387, ENSG00000233493, ENSG00000215908, ENSG00000234741, ENSG00000224066, ENSG00000215414, ENSG00000253683, ENSG00000179085"
$ws.Range("I15").Value = "SNHG5, SNHG8, RPL13AP5, RN7SL2, RP11-533E19.7, ZNF252P-AS1, FOXN3-AS1, TMEM191A, C6orf48, TMEM238, CROCCP2, GAS5, RP4-622L5.7, PSMA6P1, CTB-79E8.3, DPM3"
$ws.Range("C16").Value = "SUBSET"
$ws.Range("F16").Value = 0.902564102564102
$ws.Range("G16").Value = 0.9729783037475344
$ws.Range("H16").Value = "ENSG00000203875, ENSG00000269893, ENSG00000236552, ENSG00000274012, ENSG00000278771, ENSG00000272906, ENSG00000255559, ENSG00000258920, ENSG00000225864, ENSG00000226287, ENSG00000204387, ENSG00000233493, ENSG00000215908, ENSG00000234741, ENSG00000224066, ENSG00000253683, ENSG00000179085"
$ws.Range("I16").Value = "SNHG5, SNHG8, RPL13AP5, RN7SL2, Metazoa_SRP, RP11-533E19.7, ZNF252P-AS1, FOXN3-AS1, HCG4P11, TMEM191A, C6orf48, TMEM238, CROCCP2, GAS5, RP4-622L5.7, CTB-79E8.3, DPM3"
$ws.Range("C17").Value = "SUBSET"
$ws.Range("F17").Value = 0.8897435897435892
$ws.Range("G17").Value = 0.9705456936226163
$ws.Range("H17").Value = "ENSG00000203875, ENSG00000269893, ENSG00000236552, ENSG00000274012, ENSG00000278771, ENSG00000272906, ENSG00000255559, ENSG00000258920, ENSG00000225864, ENSG00000226287, ENSG00000204387, ENSG00000233493, ENSG00000215908, ENSG00000234741, ENSG00000224066, ENSG00000215414, ENSG00000230979, ENSG00000179085"
$ws.Range("I17").Value = "SNHG5, SNHG8, RPL13AP5, RN7SL2, Metazoa_SRP, RP11-533E19.7, ZNF252P-AS1, FOXN3-AS1, HCG4P11, TMEM191A, C6orf48, TMEM238, CROCCP2, GAS5, RP4-622L5.7, PSMA6P1, AC079250.1, DPM3"
$ws.Range("C18").Value = "SUBSET"
$ws.Range("F18").Value = 0.8974358974358971
$ws.Range("G18").Value = 0.9722550953320178
$ws.Range("H18").Value = "ENSG00000203875, ENSG00000269893, ENSG00000236552, ENSG00000274012, ENSG00000278771, ENSG00000272906, ENSG00000255559, ENSG00000258920, ENSG00000225864, ENSG00000226287, ENSG00000204387, ENSG00000233493, ENSG00000215908, ENSG00000234741, ENSG00000215414, ENSG00000230979, ENSG00000253683, ENSG00000179085, ENSG00000226085"
$ws.Range("I18").Value = "SNHG5, SNHG8, RPL13AP5, RN7SL2, Metazoa_SRP, RP11-533E19.7, ZNF252P-AS1, FOXN3-AS1, HCG4P11, TMEM191A, C6orf48, TMEM238, CROCCP2, GAS5, PSMA6P1, AC079250.1, CTB-79E8.3, DPM3, UQCRFS1P1"
$ws.Range("C19").Value = "SUBSET"
$ws.Range("F19").Value = 0.8974358974358971
$ws.Range("G19").Value = 0.9725180802103874
